# ----------------------------------------------------------------------------
# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
#
# Every Price/Volume cell in this sheet is stored as TEXT (e.g. "58.621.71",
# "  +0.27%  ") rather than as a number - so each new value below is written as
# a string. Plain assignment is enough for values Excel cannot parse as a
# number (multi-dot prices, percent strings with padding, URLs, names). For
# values that DO look like a plain number (e.g. "528.91", "1.00"), a leading
# apostrophe is used to force text, and the cell's Style is then reset back to
# 'Normal' so no stray number-format is left behind - the cell ends up exactly
# as before: plain text, default style.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$cellRef, [string]$text) {
    if ($text -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        # Looks like a plain number to Excel - force text with a quote prefix,
        # then drop back to the default style so no format override remains.
        $ws.Range($cellRef).Value = "'" + $text
        $ws.Range($cellRef).Style = "Normal"
    } else {
        $ws.Range($cellRef).Value = $text
    }
}

# Row 2: Bitcoin
Set-TextValue "D2" '58.621.71'
Set-TextValue "E2" '  +0.27%  '

# Row 3: Ethereum
Set-TextValue "D3" '3.150.74'
Set-TextValue "E3" '  +0.06%  '

# Row 4: TetherUSD
Set-TextValue "E4" '  +0.00%  '

# Row 5: BNB
Set-TextValue "D5" '528.91'
Set-TextValue "E5" '  -1.50%  '

# Row 6: Solana
Set-TextValue "D6" '139.35'
Set-TextValue "E6" '  -0.69%  '

# Row 7: USDC
Set-TextValue "E7" '  +0.11%  '

# Row 8: XRP
Set-TextValue "D8" '0.537'
Set-TextValue "E8" '  +14.37%  '

# Row 9: Toncoin
Set-TextValue "E9" '  -0.24%  '

# Row 10: Cardano
Set-TextValue "E10" '  +4.98%  '

# Row 11: Dogecoin
Set-TextValue "E11" '  +2.79%  '

# Row 12: TRON
Set-TextValue "E12" '  +3.26%  '

# Row 13: Wrapped liquid staked Ether 2.0
Set-TextValue "D13" '3.694.15'
Set-TextValue "E13" '  +0.25%  '

# Row 14: Avalanche
Set-TextValue "D14" '25.76'
Set-TextValue "E14" '  -0.45%  '

# Row 15: Shiba Inu
Set-TextValue "E15" '  +3.46%  '

# Row 16: Wrapped BTC
Set-TextValue "D16" '58.670.84'
Set-TextValue "E16" '  +0.25%  '

# Row 17: Polkadot
Set-TextValue "E17" '  +2.61%  '

# Row 18: Wrapped Ether
Set-TextValue "D18" '3.164.57'
Set-TextValue "E18" '  +0.61%  '

# Row 19: Chainlink
Set-TextValue "E19" '  +1.30%  '

# Row 20: Uniswap
Set-TextValue "D20" '8.13'
Set-TextValue "E20" '  -1.28%  '

# Row 21: Bitcoin Cash
Set-TextValue "D21" '374.07'
Set-TextValue "E21" '  +3.31%  '

# Row 22: LEO
Set-TextValue "D22" '5.79'
Set-TextValue "E22" '  +2.02%  '

# Row 23: Dai
Set-TextValue "E23" '  +0.25%  '

# Row 24: Polygon
Set-TextValue "D24" '0.529'
Set-TextValue "E24" '  +3.91%  '

# Row 25: Litecoin
Set-TextValue "D25" '69.48'
Set-TextValue "E25" '  +0.35%  '

# Row 27: Binance-Peg BSC-USD
Set-TextValue "D27" '1.00'
Set-TextValue "E27" '  +0.01%  '

# Row 28: Internet Computer (DFINITY)
Set-TextValue "D28" '8.29'
Set-TextValue "E28" '  +12.15%  '

# Row 29: PEPE
Set-TextValue "D29" '0.0₃0863'
Set-TextValue "E29" '  -2.77%  '

# Row 30: Ethereum Classic
Set-TextValue "D30" '22.21'
Set-TextValue "E30" '  +2.88%  '

# Row 31: PancakeSwap
Set-TextValue "E31" '  -0.60%  '

# Row 32: Render Token
Set-TextValue "D32" '6.08'
Set-TextValue "E32" '  -1.82%  '

# Row 33: NEAR Protocol
Set-TextValue "E33" '  -1.02%  '

# Row 34: Fetch.AI
Set-TextValue "E34" '  -0.52%  '

# Row 35: Aptos
Set-TextValue "D35" '6.26'
Set-TextValue "E35" '  +2.09%  '

# Row 36: Monero
Set-TextValue "D36" '158.01'
Set-TextValue "E36" '  -0.61%  '

# Row 37: Immutable X
Set-TextValue "E37" '  +3.77%  '

# Row 38: Energy Swap
Set-TextValue "D38" '24.92'
Set-TextValue "E38" '  -4.82%  '

# Row 39: Stacks
Set-TextValue "E39" '  -1.76%  '

# Row 40: Hedera -> Maker (row reassigned)
Set-TextValue "B40" 'Maker'
Set-TextValue "C40" 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue "D40" '2.634.32'
Set-TextValue "E40" '  +5.27%  '

# Row 41: Maker -> Hedera (row reassigned)
Set-TextValue "B41" 'Hedera'
Set-TextValue "C41" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D41" '0.0684'
Set-TextValue "E41" '  +1.22%  '

# Row 42: Filecoin
Set-TextValue "E42" '  +4.75%  '

# Row 43: Mantle
Set-TextValue "D43" '0.719'
Set-TextValue "E43" '  +1.85%  '

# Row 44: OKB
Set-TextValue "D44" '39.02'
Set-TextValue "E44" '  +3.47%  '

# Row 45: VeChain
Set-TextValue "D45" '0.0286'
Set-TextValue "E45" '  +5.85%  '

# Row 46: First Digital USD
Set-TextValue "E46" '  -0.02%  '

# Row 47: Renzo Restaked ETH
Set-TextValue "D47" '3.191.18'
Set-TextValue "E47" '  +0.06%  '

# Row 48: Stellar
Set-TextValue "E48" '  +13.00%  '

# Row 49: Cosmos
Set-TextValue "E49" '  +1.71%  '

# Row 50: ONDO
Set-TextValue "E50" '  -2.52%  '

# Row 51: Injective Protocol
Set-TextValue "D51" '19.98'
Set-TextValue "E51" '  -0.34%  '
